$d = $word.ActiveDocument

# Update the date heading (first paragraph, outside the table)
$d.Paragraphs.Item(1).Range.Text = "2023-12-24 Sunday"

# Update the math-problem table cells (5 columns x 20 rows)
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "52-25="
$t.Cell(1,2).Range.Text = "93-56="
$t.Cell(1,3).Range.Text = "90-47="
$t.Cell(1,4).Range.Text = "83-28="
$t.Cell(1,5).Range.Text = "28+43="
$t.Cell(2,1).Range.Text = "79+5="
$t.Cell(2,2).Range.Text = "35+27="
$t.Cell(2,3).Range.Text = "80-44="
$t.Cell(2,4).Range.Text = "96-69="
$t.Cell(2,5).Range.Text = "77+7="
$t.Cell(3,1).Range.Text = "18+59="
$t.Cell(3,2).Range.Text = "6+29="
$t.Cell(3,3).Range.Text = "55+38="
$t.Cell(3,4).Range.Text = "60-37="
$t.Cell(3,5).Range.Text = "36+18="
$t.Cell(4,1).Range.Text = "13+58="
$t.Cell(4,2).Range.Text = "54-39="
$t.Cell(4,3).Range.Text = "22-3="
$t.Cell(4,4).Range.Text = "65-26="
$t.Cell(4,5).Range.Text = "55-27="
$t.Cell(5,1).Range.Text = "51-4="
$t.Cell(5,2).Range.Text = "65-56="
$t.Cell(5,3).Range.Text = "81-72="
$t.Cell(5,4).Range.Text = "16+78="
$t.Cell(5,5).Range.Text = "18+77="
$t.Cell(6,1).Range.Text = "6+47="
$t.Cell(6,2).Range.Text = "91-33="
$t.Cell(6,3).Range.Text = "81-44="
$t.Cell(6,4).Range.Text = "8+46="
$t.Cell(6,5).Range.Text = "26+9="
$t.Cell(7,1).Range.Text = "23+38="
$t.Cell(7,2).Range.Text = "89+5="
$t.Cell(7,3).Range.Text = "91-28="
$t.Cell(7,4).Range.Text = "83-38="
$t.Cell(7,5).Range.Text = "4+59="
$t.Cell(8,1).Range.Text = "8+34="
$t.Cell(8,2).Range.Text = "63-9="
$t.Cell(8,3).Range.Text = "44-39="
$t.Cell(8,4).Range.Text = "46-19="
$t.Cell(8,5).Range.Text = "18+65="
$t.Cell(9,1).Range.Text = "62-7="
$t.Cell(9,2).Range.Text = "10-9="
$t.Cell(9,3).Range.Text = "81-78="
$t.Cell(9,4).Range.Text = "39+3="
$t.Cell(9,5).Range.Text = "10-9="
$t.Cell(10,1).Range.Text = "17+49="
$t.Cell(10,2).Range.Text = "64-45="
$t.Cell(10,3).Range.Text = "17+48="
$t.Cell(10,4).Range.Text = "23+49="
$t.Cell(10,5).Range.Text = "36+49="
$t.Cell(11,1).Range.Text = "49+49="
$t.Cell(11,2).Range.Text = "49+46="
$t.Cell(11,3).Range.Text = "85+8="
$t.Cell(11,4).Range.Text = "18+37="
$t.Cell(11,5).Range.Text = "94-7="
$t.Cell(12,1).Range.Text = "81-35="
$t.Cell(12,2).Range.Text = "60-14="
$t.Cell(12,3).Range.Text = "96-58="
$t.Cell(12,4).Range.Text = "64-35="
$t.Cell(12,5).Range.Text = "22-8="
$t.Cell(13,1).Range.Text = "43-29="
$t.Cell(13,2).Range.Text = "17+9="
$t.Cell(13,3).Range.Text = "94-36="
$t.Cell(13,4).Range.Text = "35-18="
$t.Cell(13,5).Range.Text = "93-19="
$t.Cell(14,2).Range.Text = "14+78="
$t.Cell(14,3).Range.Text = "76-29="
$t.Cell(14,4).Range.Text = "19+48="
$t.Cell(14,5).Range.Text = "51-19="
$t.Cell(15,1).Range.Text = "8+84="
$t.Cell(15,2).Range.Text = "25+57="
$t.Cell(15,3).Range.Text = "6+75="
$t.Cell(15,4).Range.Text = "33+58="
$t.Cell(15,5).Range.Text = "19+64="
$t.Cell(16,1).Range.Text = "13+68="
$t.Cell(16,2).Range.Text = "63-19="
$t.Cell(16,3).Range.Text = "92-78="
$t.Cell(16,4).Range.Text = "41-2="
$t.Cell(16,5).Range.Text = "31-23="
$t.Cell(17,1).Range.Text = "9+78="
$t.Cell(17,2).Range.Text = "65+28="
$t.Cell(17,3).Range.Text = "42-14="
$t.Cell(17,4).Range.Text = "79+17="
$t.Cell(17,5).Range.Text = "84-29="
$t.Cell(18,1).Range.Text = "18+34="
$t.Cell(18,2).Range.Text = "28+16="
$t.Cell(18,3).Range.Text = "41-3="
$t.Cell(18,4).Range.Text = "84-79="
$t.Cell(18,5).Range.Text = "48+5="
$t.Cell(19,1).Range.Text = "35+56="
$t.Cell(19,2).Range.Text = "28-9="
$t.Cell(19,3).Range.Text = "50-28="
$t.Cell(19,4).Range.Text = "48+13="
$t.Cell(19,5).Range.Text = "31-7="
$t.Cell(20,1).Range.Text = "12-3="
$t.Cell(20,2).Range.Text = "18+48="
$t.Cell(20,3).Range.Text = "9+49="
$t.Cell(20,4).Range.Text = "82-27="
$t.Cell(20,5).Range.Text = "46+25="
